{"js": "// Rewrites the three paragraphs of tscp.docx to the \"ordinal\" wording\n// used by the TSCP/BAILS smart-tag test, and maps the smart tag's\n// attributes from the old BAF names to the new BAILS names.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst p1 = paragraphs.items[0];\nconst p2 = paragraphs.items[1];\nconst p3 = paragraphs.items[2];\n\n// --- Paragraph 1: \"before\" -> \"1st paragraph, non-business.\" -------------\np1.clear();\nlet r = p1.insertText(\"1\", Word.InsertLocation.start);\nawait context.sync();\nconst sup1 = r.insertText(\"st\", Word.InsertLocation.after);\nawait context.sync();\nsup1.insertText(\" paragraph, non-business.\", Word.InsertLocation.after);\nawait context.sync();\n// Only the ordinal suffix (\"st\") is superscript.\nsup1.font.superscript = true;\np1.style = \"Normal\";\nawait context.sync();\n\n// --- Paragraph 3: \"after\" -> \"3rd paragraph, non-business.\" --------------\np3.clear();\nr = p3.insertText(\"3\", Word.InsertLocation.start);\nawait context.sync();\nconst sup3 = r.insertText(\"rd\", Word.InsertLocation.after);\nawait context.sync();\nsup3.insertText(\" paragraph, non-business.\", Word.InsertLocation.after);\nawait context.sync();\nsup3.font.superscript = true;\np3.style = \"Normal\";\nawait context.sync();\n\n// --- Paragraph 2: holds the RDF smart tag ---------------------------------\n// Rewrite it wholesale via WordprocessingML so we can, in one go: drop the\n// stale \"_GoBack\" bookmark, remap the smart tag attribute names from the\n// BAF namespace to BAILS, and rebuild the \"2nd paragraph, confidential.\"\n// run sequence (with its superscript \"nd\").\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"Normal\"/>\n            </w:pPr>\n            <w:smartTag w:uri=\"http://www.w3.org/1999/02/22-rdf-syntax-ns#\" w:element=\"RDF\">\n              <w:smartTagPr>\n                <w:attr w:name=\"urn:bails:ExportControl:Authorization:StartValidity\" w:val=\"2015-11-27\"/>\n                <w:attr w:name=\"urn:bails:ExportControl:BusinessAuthorization:Identifier\" w:val=\"urn:example:tscp:1\"/>\n                <w:attr w:name=\"urn:bails:ExportControl:BusinessAuthorizationCategory:Identifier\" w:val=\"urn:example:tscp:1:confidential\"/>\n              </w:smartTagPr>\n            </w:smartTag>\n            <w:r>\n              <w:t>2</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:vertAlign w:val=\"superscript\"/>\n              </w:rPr>\n              <w:t>nd</w:t>\n            </w:r>\n            <w:r>\n              <w:t xml:space=\"preserve\"> paragraph, confidential.</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\np2.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Rewrites the three paragraphs of tscp.docx to the \"ordinal\" wording\n# used by the TSCP/BAILS smart-tag test, and maps the smart tag's\n# attributes from the old BAF names to the new BAILS names.\n\n$d = $word.ActiveDocument\n\nfunction Set-OrdinalParagraph($para, [string]$num, [string]$suffix, [string]$rest) {\n    # Put the leading digit into the paragraph first.\n    $para.Range.Text = $num\n    $pos = $para.Range.Start + $num.Length\n\n    # Insert the ordinal suffix (\"st\"/\"nd\"/\"rd\") right after the digit,\n    # then mark just that piece as superscript.\n    $r = $d.Range($pos, $pos)\n    $null = $r.InsertAfter($suffix)\n    $supRange = $d.Range($pos, $pos + $suffix.Length)\n    $supRange.Font.Superscript = $true\n\n    # Append the remaining (non-superscript) text after the suffix.\n    $afterPos = $pos + $suffix.Length\n    $r2 = $d.Range($afterPos, $afterPos)\n    $null = $r2.InsertAfter($rest)\n}\n\n# Paragraph 1: \"before\" -> \"1st paragraph, non-business.\"\nSet-OrdinalParagraph $d.Paragraphs(1) \"1\" \"st\" \" paragraph, non-business.\"\n$d.Paragraphs(1).Style = \"Normal\"\n\n# Paragraph 3: \"after\" -> \"3rd paragraph, non-business.\"\nSet-OrdinalParagraph $d.Paragraphs(3) \"3\" \"rd\" \" paragraph, non-business.\"\n$d.Paragraphs(3).Style = \"Normal\"\n\n# Paragraph 2 holds the RDF smart tag. Rewrite its contents wholesale via\n# WordprocessingML so we can: drop the stale \"_GoBack\" bookmark, remap the\n# smart tag attribute names from the BAF namespace to BAILS, and rebuild\n# the \"2nd paragraph, confidential.\" run sequence with its superscript.\n$p2 = $d.Paragraphs(2)\n$p2xml = @'\n<?xml version=\"1.0\" standalone=\"yes\"?>\n<?mso-application progid=\"Word.Document\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"Normal\"/>\n            </w:pPr>\n            <w:smartTag w:uri=\"http://www.w3.org/1999/02/22-rdf-syntax-ns#\" w:element=\"RDF\">\n              <w:smartTagPr>\n                <w:attr w:name=\"urn:bails:ExportControl:Authorization:StartValidity\" w:val=\"2015-11-27\"/>\n                <w:attr w:name=\"urn:bails:ExportControl:BusinessAuthorization:Identifier\" w:val=\"urn:example:tscp:1\"/>\n                <w:attr w:name=\"urn:bails:ExportControl:BusinessAuthorizationCategory:Identifier\" w:val=\"urn:example:tscp:1:confidential\"/>\n              </w:smartTagPr>\n            </w:smartTag>\n            <w:r>\n              <w:t>2</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:vertAlign w:val=\"superscript\"/>\n              </w:rPr>\n              <w:t>nd</w:t>\n            </w:r>\n            <w:r>\n              <w:t xml:space=\"preserve\"> paragraph, confidential.</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n$null = $p2.Range.InsertXML($p2xml)\n"}
